$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto snapshot refresh: update the Price (D) and
# Volume(1h) (E) columns for every listed coin row.
#
# The Price column stores plain text (e.g. "315.20", "0.000008902")
# even though most values look numeric. Excel auto-converts a bare
# numeric-looking string assigned via .Value into a real number,
# which silently drops significant trailing zeros (315.20 -> 315.2)
# or renders tiny values in scientific notation. Setting the cell to
# Text format first keeps the exact original string representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.388.82"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.49"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -3.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.20"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4294"
$ws.Range("E7").Value = "  -2.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3699"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07251"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8666"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.15"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.48"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.678"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.359"
$ws.Range("E14").Value = "  -3.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07098"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.70"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008902"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.409.06"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.169"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.061.05"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.014"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.32"
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.48"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.141"
$ws.Range("E28").Value = "  +6.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.292"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.46"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08862"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.204"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7670"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.501"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.893"
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.122"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01961"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05278"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.164"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.886"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1676"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5071"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.667"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.58"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.32"
$ws.Range("E46").Value = "  -4.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4729"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06432"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.669"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.824"
$ws.Range("E51").Value = "  -3.60%  "
